$d = $word.ActiveDocument

# --- 1. Move the "_GoBack" bookmark from right after "MP73010" (inside
#        the title paragraph) to the end of the
#        ">>> ... stuff after this line >>>" paragraph, right before its
#        paragraph mark. ---
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

$targetPara = $d.Paragraphs.Item(4)
$pos = $targetPara.Range.End - 1

# A collapsed (zero-length) Range placed exactly one character before a
# paragraph mark cannot be bookmarked directly here, so nudge it into
# place: insert a throwaway character, bookmark the one-character range
# around it, then delete the throwaway character. The bookmark collapses
# to the desired position and survives the deletion.
$d.Range($pos, $pos).InsertAfter("X")
$markRange = $d.Range($pos, $pos + 1)
$d.Bookmarks.Add("_GoBack", $markRange)
$d.Range($pos, $pos + 1).Delete()

# --- 2. Turn the last (empty) paragraph into the new DVCS paragraph,
#        keeping the paragraph before it empty. ---
$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)

$newParaXml = @'
<?xml version="1.0" standalone="yes"?>
<?mso-application progid="Word.Document"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
<w:p>
<w:pPr><w:rPr><w:color w:val="0D0D0D" w:themeColor="text1" w:themeTint="F2"/></w:rPr></w:pPr>
<w:r>
<w:rPr><w:rFonts w:cs="Arial"/><w:color w:val="0D0D0D" w:themeColor="text1" w:themeTint="F2"/><w:shd w:val="clear" w:color="auto" w:fill="FCFCFA"/></w:rPr>
<w:t xml:space="preserve">Distributed Version Control Systems (DVCSs) step in. In a DVCS (such as Git, Mercurial, Bazaar or </w:t>
</w:r>
<w:proofErr w:type="spellStart"/>
<w:r>
<w:rPr><w:rFonts w:cs="Arial"/><w:color w:val="0D0D0D" w:themeColor="text1" w:themeTint="F2"/><w:shd w:val="clear" w:color="auto" w:fill="FCFCFA"/></w:rPr>
<w:t>Darcs</w:t>
</w:r>
<w:proofErr w:type="spellEnd"/>
<w:r>
<w:rPr><w:rFonts w:cs="Arial"/><w:color w:val="0D0D0D" w:themeColor="text1" w:themeTint="F2"/><w:shd w:val="clear" w:color="auto" w:fill="FCFCFA"/></w:rPr>
<w:t>), clients don&#8217;t just check out the latest snapshot of the files; rather, they fully mirror the repository, including its full history. Thus, if any server dies, and these systems were collaborating via that server, any of the client repositories can be copied back up to the server to restore it. Every clone is really a full backup of all the data</w:t>
</w:r>
<w:r>
<w:rPr><w:rFonts w:cs="Arial"/><w:color w:val="0D0D0D" w:themeColor="text1" w:themeTint="F2"/><w:shd w:val="clear" w:color="auto" w:fill="FCFCFA"/></w:rPr>
<w:t>.</w:t>
</w:r>
</w:p>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
'@

[void]$lastPara.Range.InsertXML($newParaXml)

# InsertXML inserted the new paragraph *before* the previously-empty
# last paragraph, leaving that empty paragraph trailing behind it. Merge
# that now-redundant trailing empty paragraph away so the document ends
# with exactly one empty paragraph followed by the new content paragraph.
$finalPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$mergeRange = $d.Range($finalPara.Range.Start - 1, $finalPara.Range.End)
$mergeRange.Delete()
